# Sprint 10 and PB update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint")

# Asigned person updates (column D)
$ws.Range("D7").Value = "Đuro Nenadović"
$ws.Range("D10").Value = "Vanja Cvetković"

# Task status / progress updates (row 9)
$ws.Range("F9").Value = "Done"
$ws.Range("K9").Value = 5

# Update the sheet view: scroll so column B is the left-most visible column,
# and move the active selection to K9
$ws.Activate()
$ws.Range("K9").Select()
$excel.ActiveWindow.ScrollColumn = 2

$wb.Save()
